$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.684.02"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").Value = "2.306.96"
$ws.Range("E3").Value = "  -4.16%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "131.87"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").Value = "2.305.20"
$ws.Range("E9").Value = "  -4.17%  "
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("D11").Value = "5.50"
$ws.Range("E11").Value = "  -2.30%  "
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("E13").Value = "  -4.63%  "
$ws.Range("D14").Value = "23.90"
$ws.Range("E14").Value = "  -2.91%  "
$ws.Range("D15").Value = "2.718.61"
$ws.Range("E15").Value = "  -4.19%  "
$ws.Range("D16").Value = "58.704.06"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "2.315.15"
$ws.Range("E18").Value = "  -3.75%  "
$ws.Range("D19").Value = "10.68"
$ws.Range("E19").Value = "  -4.38%  "
$ws.Range("D20").Value = "4.31"
$ws.Range("E20").Value = "  -4.05%  "
$ws.Range("D21").Value = "314.57"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  -4.49%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "63.32"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E25").Value = "  -6.76%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "8.09"
$ws.Range("E27").Value = "  -5.91%  "
$ws.Range("E28").Value = "  -6.04%  "
$ws.Range("D29").Value = "1.76"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "168.49"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "0.0₃0725"
$ws.Range("E31").Value = "  -5.30%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -5.61%  "
$ws.Range("E34").Value = "  -4.75%  "
$ws.Range("D36").Value = "17.78"
$ws.Range("E36").Value = "  -3.31%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "1.25"
$ws.Range("E38").Value = "  -5.19%  "
$ws.Range("E39").Value = "  -5.11%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -4.98%  "
$ws.Range("D42").Value = "297.39"
$ws.Range("E42").Value = "  -7.86%  "
$ws.Range("D43").Value = "141.63"
$ws.Range("E43").Value = "  -3.74%  "
$ws.Range("E44").Value = "  -4.00%  "
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "18.44"
$ws.Range("E48").Value = "  -7.02%  "
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("D50").Value = "16.64"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("E51").Value = "  -0.43%  "
